$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "Resolución Primigenia"
$ws.Range("D1").Value = "Resolución Hija"
